# Natmi following Dr Hou advice: update the Il27-Il27ra ligand-receptor
# pairs table (Sheet1) with recomputed statistics across 3 cell
# replicates, and add the new "M1" and "sCs" target-cluster rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Il27"
$ws.Range("C2").Value = "Il27ra"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9134243333333334
$ws.Range("H2").Value = 2.740273
$ws.Range("I2").Value = 0.4164467621298553
$ws.Range("J2").Value = 0.4164467621298553
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.977152666666667
$ws.Range("N2").Value = 5.931458
$ws.Range("O2").Value = 0.5296093122610759
$ws.Range("P2").Value = 0.5296093122610759
$ws.Range("Q2").Value = 1.805979356448222
$ws.Range("R2").Value = 16.253814208034
$ws.Range("S2").Value = 0.2205540832849445
$ws.Range("T2").Value = 0.2205540832849445

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Il27"
$ws.Range("C3").Value = "Il27ra"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9134243333333334
$ws.Range("H3").Value = 2.740273
$ws.Range("I3").Value = 0.4164467621298553
$ws.Range("J3").Value = 0.4164467621298553
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.561521333333333
$ws.Range("N3").Value = 4.684564
$ws.Range("O3").Value = 0.4182763695339315
$ws.Range("P3").Value = 0.4182763695339315
$ws.Range("Q3").Value = 1.426331582885778
$ws.Range("R3").Value = 12.836984245972
$ws.Range("S3").Value = 0.1741898397678366
$ws.Range("T3").Value = 0.1741898397678366

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Il27"
$ws.Range("C4").Value = "Il27ra"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9134243333333334
$ws.Range("H4").Value = 2.740273
$ws.Range("I4").Value = 0.4164467621298553
$ws.Range("J4").Value = 0.4164467621298553
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.04091566666666666
$ws.Range("N4").Value = 0.122747
$ws.Range("O4").Value = 0.01095986083895566
$ws.Range("P4").Value = 0.01095986083895566
$ws.Range("Q4").Value = 0.03737336554788889
$ws.Range("R4").Value = 0.336360289931
$ws.Range("S4").Value = 0.004564198559776885
$ws.Range("T4").Value = 0.004564198559776885

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Il27"
$ws.Range("C5").Value = "Il27ra"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9134243333333334
$ws.Range("H5").Value = 2.740273
$ws.Range("I5").Value = 0.4164467621298553
$ws.Range("J5").Value = 0.4164467621298553
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1059233333333333
$ws.Range("N5").Value = 0.31777
$ws.Range("O5").Value = 0.02837311688916993
$ws.Range("P5").Value = 0.02837311688916993
$ws.Range("Q5").Value = 0.09675295013444445
$ws.Range("R5").Value = 0.8707765512100001
$ws.Range("S5").Value = 0.01181589266002673
$ws.Range("T5").Value = 0.01181589266002673

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Il27"
$ws.Range("C6").Value = "Il27ra"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9134243333333334
$ws.Range("H6").Value = 2.740273
$ws.Range("I6").Value = 0.4164467621298553
$ws.Range("J6").Value = 0.4164467621298553
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.04771566666666666
$ws.Range("N6").Value = 0.143147
$ws.Range("O6").Value = 0.01278134047686694
$ws.Range("P6").Value = 0.01278134047686694
$ws.Range("Q6").Value = 0.04358465101455555
$ws.Range("R6").Value = 0.392261859131
$ws.Range("S6").Value = 0.005322747857270498
$ws.Range("T6").Value = 0.005322747857270497

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Il27"
$ws.Range("C7").Value = "Il27ra"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.279951666666667
$ws.Range("H7").Value = 3.839855
$ws.Range("I7").Value = 0.5835532378701448
$ws.Range("J7").Value = 0.5835532378701448
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.977152666666667
$ws.Range("N7").Value = 5.931458
$ws.Range("O7").Value = 0.5296093122610759
$ws.Range("P7").Value = 0.5296093122610759
$ws.Range("Q7").Value = 2.530659850954444
$ws.Range("R7").Value = 22.77593865859
$ws.Range("S7").Value = 0.3090552289761314
$ws.Range("T7").Value = 0.3090552289761314

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Il27"
$ws.Range("C8").Value = "Il27ra"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.279951666666667
$ws.Range("H8").Value = 3.839855
$ws.Range("I8").Value = 0.5835532378701448
$ws.Range("J8").Value = 0.5835532378701448
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.561521333333333
$ws.Range("N8").Value = 4.684564
$ws.Range("O8").Value = 0.4182763695339315
$ws.Range("P8").Value = 0.4182763695339315
$ws.Range("Q8").Value = 1.998671833135556
$ws.Range("R8").Value = 17.98804649822
$ws.Range("S8").Value = 0.2440865297660949
$ws.Range("T8").Value = 0.2440865297660949

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Il27"
$ws.Range("C9").Value = "Il27ra"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.279951666666667
$ws.Range("H9").Value = 3.839855
$ws.Range("I9").Value = 0.5835532378701448
$ws.Range("J9").Value = 0.5835532378701448
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.04091566666666666
$ws.Range("N9").Value = 0.122747
$ws.Range("O9").Value = 0.01095986083895566
$ws.Range("P9").Value = 0.01095986083895566
$ws.Range("Q9").Value = 0.05237007574277777
$ws.Range("R9").Value = 0.471330681685
$ws.Range("S9").Value = 0.006395662279178778
$ws.Range("T9").Value = 0.006395662279178778

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Il27"
$ws.Range("C10").Value = "Il27ra"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.279951666666667
$ws.Range("H10").Value = 3.839855
$ws.Range("I10").Value = 0.5835532378701448
$ws.Range("J10").Value = 0.5835532378701448
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.1059233333333333
$ws.Range("N10").Value = 0.31777
$ws.Range("O10").Value = 0.02837311688916993
$ws.Range("P10").Value = 0.02837311688916993
$ws.Range("Q10").Value = 0.1355767470388889
$ws.Range("R10").Value = 1.22019072335
$ws.Range("S10").Value = 0.0165572242291432
$ws.Range("T10").Value = 0.0165572242291432

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Il27"
$ws.Range("C11").Value = "Il27ra"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.279951666666667
$ws.Range("H11").Value = 3.839855
$ws.Range("I11").Value = 0.5835532378701448
$ws.Range("J11").Value = 0.5835532378701448
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.04771566666666666
$ws.Range("N11").Value = 0.143147
$ws.Range("O11").Value = 0.01278134047686694
$ws.Range("P11").Value = 0.01278134047686694
$ws.Range("Q11").Value = 0.0610737470761111
$ws.Range("R11").Value = 0.5496637236849999
$ws.Range("S11").Value = 0.007458592619596443
$ws.Range("T11").Value = 0.007458592619596443
